$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").Value = $null

$ws.Range("H38").Value = 492.57144
$ws.Range("I38").Value = 112
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 336
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = 36
$ws.Range("N38").Value = -3744

$ws.Range("H39").Value = 369.42856
$ws.Range("J39").Value = 512.25
$ws.Range("L39").Value = 1536.75
$ws.Range("N39").Value = -2128.75

$ws.Range("H42").Value = 266.75
$ws.Range("I42").Value = 60
$ws.Range("K42").Value = 180
$ws.Range("M42").Value = 50

$ws.Range("H129").Value = 684
$ws.Range("I129").Value = 440
$ws.Range("J129").Value = 1050
$ws.Range("K129").Value = 1320
$ws.Range("L129").Value = 3150
$ws.Range("M129").Value = 3680
$ws.Range("N129").Value = -13150

$ws.Range("H132").Value = 4574.36
$ws.Range("I132").Value = 5318.8945
$ws.Range("J132").Value = 2216.6667
$ws.Range("K132").Value = 15956.6835
$ws.Range("L132").Value = 6650.000100000001
$ws.Range("M132").Value = -13426.6835
$ws.Range("N132").Value = -11710.0001

$ws.Range("H137").Value = 23257406
$ws.Range("I137").Value = 1034.2963
$ws.Range("J137").Value = 62502536
$ws.Range("K137").Value = 3102.8889
$ws.Range("L137").Value = 187507608
$ws.Range("M137").Value = -552.8888999999999
$ws.Range("N137").Value = -187512708

$ws.Range("H138").Value = 2366.3447
$ws.Range("I138").Value = 2788.0688
$ws.Range("J138").Value = 2155.4827
$ws.Range("K138").Value = 8364.206399999999
$ws.Range("L138").Value = 6466.4481
$ws.Range("M138").Value = -3224.206399999999
$ws.Range("N138").Value = -16746.4481

$ws.Range("H141").Value = 1801.3077
$ws.Range("I141").Value = 831.7
$ws.Range("K141").Value = 2495.1
$ws.Range("M141").Value = 2684.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5485.9443
$ws.Range("I61").Value = 6467.9165
$ws.Range("J61").Value = 3522
$ws.Range("K61").Value = 6467.9165
$ws.Range("L61").Value = 3522
$ws.Range("M61").Value = -6255.9165
$ws.Range("N61").Value = -3946

$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = $null
$ws.Range("M122").Value = $null
$ws.Range("N122").Value = 0

$ws.Range("H132").Value = 57770.61
$ws.Range("I132").Value = 2176.7778
$ws.Range("K132").Value = 6530.3334
$ws.Range("M132").Value = -4000.3334

$ws.Range("H136").Value = 5485.9443
$ws.Range("I136").Value = 6467.9165
$ws.Range("J136").Value = 3522
$ws.Range("K136").Value = 19403.7495
$ws.Range("L136").Value = 10566
$ws.Range("M136").Value = -16853.7495
$ws.Range("N136").Value = -15666

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 468.57144
$ws.Range("I94").Value = 463.33334
$ws.Range("K94").Value = 463.33334
$ws.Range("M94").Value = -12.33334000000002

$ws.Range("H134").Value = 80775.57000000001
$ws.Range("I134").Value = 158953.72
$ws.Range("J134").Value = 2597.4285
$ws.Range("K134").Value = 476861.16
$ws.Range("L134").Value = 7792.2855
$ws.Range("M134").Value = -474326.16
$ws.Range("N134").Value = -12862.2855

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4568.8125
$ws.Range("I86").Value = 4100.091
$ws.Range("J86").Value = 5600
$ws.Range("K86").Value = 4100.091
$ws.Range("L86").Value = 5600
$ws.Range("M86").Value = -2977.091
$ws.Range("N86").Value = -7846

$ws.Range("H89").Value = 4568.8125
$ws.Range("I89").Value = 4100.091
$ws.Range("J89").Value = 5600
$ws.Range("K89").Value = 20500.455
$ws.Range("L89").Value = 28000
$ws.Range("M89").Value = -14884.455
$ws.Range("N89").Value = -39232

$ws.Range("H99").Value = 64032.582
$ws.Range("I99").Value = 34862.8
$ws.Range("J99").Value = 112648.89
$ws.Range("K99").Value = 34862.8
$ws.Range("L99").Value = 112648.89
$ws.Range("M99").Value = -33364.8
$ws.Range("N99").Value = -115644.89

$ws.Range("H126").Value = 64032.582
$ws.Range("I126").Value = 34862.8
$ws.Range("J126").Value = 112648.89
$ws.Range("K126").Value = 104588.4
$ws.Range("L126").Value = 337946.67
$ws.Range("M126").Value = -102118.4
$ws.Range("N126").Value = -342886.67

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1001.6667
$ws.Range("I51").Value = 100
$ws.Range("J51").Value = 1452.5
$ws.Range("K51").Value = 300
$ws.Range("L51").Value = 4357.5
$ws.Range("M51").Value = 160
$ws.Range("N51").Value = -5277.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2375.375
$ws.Range("I97").Value = 1994
$ws.Range("J97").Value = 3011
$ws.Range("K97").Value = 1994
$ws.Range("L97").Value = 3011
$ws.Range("M97").Value = -1498
$ws.Range("N97").Value = -4003

$ws.Range("H122").Value = 3423.8823
$ws.Range("I122").Value = 3900.6
$ws.Range("J122").Value = 2742.8572
$ws.Range("K122").Value = 11701.8
$ws.Range("L122").Value = 8228.571599999999
$ws.Range("M122").Value = -9251.799999999999
$ws.Range("N122").Value = -13128.5716

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2088.182
$ws.Range("I122").Value = 1775
$ws.Range("K122").Value = 5325
$ws.Range("M122").Value = -2875

$ws.Range("H136").Value = 1286.826
$ws.Range("I136").Value = 963
$ws.Range("J136").Value = 2825
$ws.Range("K136").Value = 2889
$ws.Range("L136").Value = 8475
$ws.Range("M136").Value = -339
$ws.Range("N136").Value = -13575

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("L94").Value = $null
$ws.Range("N94").Value = 0

$ws.Range("H122").Value = 3502.6667
$ws.Range("I122").Value = 2146.5
$ws.Range("J122").Value = 4858.8335
$ws.Range("K122").Value = 6439.5
$ws.Range("L122").Value = 14576.5005
$ws.Range("M122").Value = -3989.5
$ws.Range("N122").Value = -19476.5005

$ws.Range("H136").Value = 5623.857
$ws.Range("I136").Value = 6508.7144
$ws.Range("J136").Value = 2969.2856
$ws.Range("K136").Value = 19526.1432
$ws.Range("L136").Value = 8907.856800000001
$ws.Range("M136").Value = -16976.1432
$ws.Range("N136").Value = -14007.8568

Write-Output "Applied Ifrit_Profits market-data refresh"